$d = $word.ActiveDocument

# 1) Merge the two runs around the old _GoBack bookmark (between "causes" and
#    ", endocrine disorders") back into one run. Re-"replacing" the exact
#    same text forces Word to re-flow/merge the run boundary, which also
#    clears out the old _GoBack bookmark that sat there.
$rng = $d.Content
$rng.Find.Execute("maternal and neonatal causes, endocrine", $true, $false, $false, $false, $false, $true, 1, $false, "maternal and neonatal causes, endocrine", 2) | Out-Null

# 2) Relocate the _GoBack bookmark: it now belongs between "...refer to th"
#    and "e ..." in the very first paragraph (this is where the author's
#    last edit landed).
$rng2 = $d.Content
$rng2.Find.Execute("Note: All page numbers refer to th", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($rng2.End, $rng2.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 3) "...Table 1; P. 10-11)." -> "...Table 1; PP. 10-11)."
$rng3 = $d.Content
$rng3.Find.Execute("Table 1; P. 10-11)", $true, $false, $false, $false, $false, $true, 1, $false, "Table 1; PP. 10-11)", 2) | Out-Null

# 4) "...Figure 3 and PP. 12)..." -> "...Figure 3 and P. 12)..."
$rng4 = $d.Content
$rng4.Find.Execute("Figure 3 and PP. 12)", $true, $false, $false, $false, $false, $true, 1, $false, "Figure 3 and P. 12)", 2) | Out-Null

# 5) "...Figure 2 (PP. 11, 12) ..." -> "...Figure 2 (PP. 11-12) ..."
$rng5 = $d.Content
$rng5.Find.Execute("PP. 11, 12)", $true, $false, $false, $false, $false, $true, 1, $false, "PP. 11-12)", 2) | Out-Null

# 6) "...Figure 5 and PP. 12)." -> "...Figure 5 and P. 12)."
$rng6 = $d.Content
$rng6.Find.Execute("Figure 5 and PP. 12)", $true, $false, $false, $false, $false, $true, 1, $false, "Figure 5 and P. 12)", 2) | Out-Null

# 7) "...as suggested (PP. 6,11)." -> "...as suggested (P. 6, P. 11)."
$rng7 = $d.Content
$rng7.Find.Execute("(PP. 6,11).", $true, $false, $false, $false, $false, $true, 1, $false, "(P. 6, P. 11).", 2) | Out-Null
